# Update NATMI ligand-receptor (Bdnf-Sort1) edge table with new TPM-derived values.
# The sending-cluster set now includes "ECs" in addition to "MuSCs", so the table
# grows from 6 data rows to 12 (rows 2-7 = ECs sender, rows 8-13 = MuSCs sender),
# and every numeric "edge" metric (M:T, plus derived specificities) is recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20

# Row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Bdnf"
$data[0,2] = "Sort1"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.08140533333333333
$data[0,7] = 0.244216
$data[0,8] = 0.1131514935296598
$data[0,9] = 0.1131514935296598
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 1.837384
$data[0,13] = 5.512152
$data[0,14] = 0.0635335947613339
$data[0,15] = 0.0635335947613339
$data[0,16] = 0.1495728569813333
$data[0,17] = 1.346155712832
$data[0,18] = 0.007188921136553099
$data[0,19] = 0.007188921136553098

# Row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Bdnf"
$data[1,2] = "Sort1"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.08140533333333333
$data[1,7] = 0.244216
$data[1,8] = 0.1131514935296598
$data[1,9] = 0.1131514935296598
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.6023626666666666
$data[1,13] = 1.807088
$data[1,14] = 0.02082867030699976
$data[1,15] = 0.02082867030699976
$data[1,16] = 0.04903553366755554
$data[1,17] = 0.4413198030079999
$data[1,18] = 0.0023567951534739
$data[1,19] = 0.0023567951534739

# Row 4: ECs -> Inflammatory-Mac
$data[2,0] = "ECs"
$data[2,1] = "Bdnf"
$data[2,2] = "Sort1"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.08140533333333333
$data[2,7] = 0.244216
$data[2,8] = 0.1131514935296598
$data[2,9] = 0.1131514935296598
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 5.528959666666668
$data[2,13] = 16.586879
$data[2,14] = 0.1911819646376369
$data[2,15] = 0.1911819646376369
$data[2,16] = 0.4500868046515556
$data[2,17] = 4.050781241864001
$data[2,18] = 0.02163252483468322
$data[2,19] = 0.02163252483468322

# Row 5: ECs -> MuSCs
$data[3,0] = "ECs"
$data[3,1] = "Bdnf"
$data[3,2] = "Sort1"
$data[3,3] = "MuSCs"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.08140533333333333
$data[3,7] = 0.244216
$data[3,8] = 0.1131514935296598
$data[3,9] = 0.1131514935296598
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 3.046454666666667
$data[3,13] = 9.139364
$data[3,14] = 0.1053411895666744
$data[3,15] = 0.1053411895666744
$data[3,16] = 0.2479976576248889
$data[3,17] = 2.231978918624
$data[3,18] = 0.01191951292966023
$data[3,19] = 0.01191951292966023

# Row 6: ECs -> Neutrophils
$data[4,0] = "ECs"
$data[4,1] = "Bdnf"
$data[4,2] = "Sort1"
$data[4,3] = "Neutrophils"
$data[4,4] = 1
$data[4,5] = 0.3333333333333333
$data[4,6] = 0.08140533333333333
$data[4,7] = 0.244216
$data[4,8] = 0.1131514935296598
$data[4,9] = 0.1131514935296598
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 9.021246333333332
$data[4,13] = 27.063739
$data[4,14] = 0.3119392618985303
$data[4,15] = 0.3119392618985303
$data[4,16] = 0.734377564847111
$data[4,17] = 6.609398083623999
$data[4,18] = 0.0352963933743584
$data[4,19] = 0.0352963933743584

# Row 7: ECs -> Resolving-Mac
$data[5,0] = "ECs"
$data[5,1] = "Bdnf"
$data[5,2] = "Sort1"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 1
$data[5,5] = 0.3333333333333333
$data[5,6] = 0.08140533333333333
$data[5,7] = 0.244216
$data[5,8] = 0.1131514935296598
$data[5,9] = 0.1131514935296598
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 8.883473666666667
$data[5,13] = 26.650421
$data[5,14] = 0.3071753188288246
$data[5,15] = 0.3071753188288246
$data[5,16] = 0.7231621349928889
$data[5,17] = 6.508459214936
$data[5,18] = 0.03475734610093092
$data[5,19] = 0.03475734610093092

# Row 8: MuSCs -> ECs
$data[6,0] = "MuSCs"
$data[6,1] = "Bdnf"
$data[6,2] = "Sort1"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.6380313333333333
$data[6,7] = 1.914094
$data[6,8] = 0.8868485064703402
$data[6,9] = 0.8868485064703401
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 1.837384
$data[6,13] = 5.512152
$data[6,14] = 0.0635335947613339
$data[6,15] = 0.0635335947613339
$data[6,16] = 1.172308563365333
$data[6,17] = 10.550777070288
$data[6,18] = 0.05634467362478079
$data[6,19] = 0.05634467362478079

# Row 9: MuSCs -> FAPs
$data[7,0] = "MuSCs"
$data[7,1] = "Bdnf"
$data[7,2] = "Sort1"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.6380313333333333
$data[7,7] = 1.914094
$data[7,8] = 0.8868485064703402
$data[7,9] = 0.8868485064703401
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.6023626666666666
$data[7,13] = 1.807088
$data[7,14] = 0.02082867030699976
$data[7,15] = 0.02082867030699976
$data[7,16] = 0.3843262553635555
$data[7,17] = 3.458936298272
$data[7,18] = 0.01847187515352586
$data[7,19] = 0.01847187515352586

# Row 10: MuSCs -> Inflammatory-Mac
$data[8,0] = "MuSCs"
$data[8,1] = "Bdnf"
$data[8,2] = "Sort1"
$data[8,3] = "Inflammatory-Mac"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.6380313333333333
$data[8,7] = 1.914094
$data[8,8] = 0.8868485064703402
$data[8,9] = 0.8868485064703401
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 5.528959666666668
$data[8,13] = 16.586879
$data[8,14] = 0.1911819646376369
$data[8,15] = 0.1911819646376369
$data[8,16] = 3.527649508069556
$data[8,17] = 31.74884557262601
$data[8,18] = 0.1695494398029537
$data[8,19] = 0.1695494398029536

# Row 11: MuSCs -> MuSCs
$data[9,0] = "MuSCs"
$data[9,1] = "Bdnf"
$data[9,2] = "Sort1"
$data[9,3] = "MuSCs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.6380313333333333
$data[9,7] = 1.914094
$data[9,8] = 0.8868485064703402
$data[9,9] = 0.8868485064703401
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 3.046454666666667
$data[9,13] = 9.139364
$data[9,14] = 0.1053411895666744
$data[9,15] = 0.1053411895666744
$data[9,16] = 1.943733532912889
$data[9,17] = 17.493601796216
$data[9,18] = 0.09342167663701421
$data[9,19] = 0.0934216766370142

# Row 12: MuSCs -> Neutrophils
$data[10,0] = "MuSCs"
$data[10,1] = "Bdnf"
$data[10,2] = "Sort1"
$data[10,3] = "Neutrophils"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 0.6380313333333333
$data[10,7] = 1.914094
$data[10,8] = 0.8868485064703402
$data[10,9] = 0.8868485064703401
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 9.021246333333332
$data[10,13] = 27.063739
$data[10,14] = 0.3119392618985303
$data[10,15] = 0.3119392618985303
$data[10,16] = 5.75583782638511
$data[10,17] = 51.80254043746599
$data[10,18] = 0.2766428685241719
$data[10,19] = 0.2766428685241719

# Row 13: MuSCs -> Resolving-Mac
$data[11,0] = "MuSCs"
$data[11,1] = "Bdnf"
$data[11,2] = "Sort1"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.6380313333333333
$data[11,7] = 1.914094
$data[11,8] = 0.8868485064703402
$data[11,9] = 0.8868485064703401
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 8.883473666666667
$data[11,13] = 26.650421
$data[11,14] = 0.3071753188288246
$data[11,15] = 0.3071753188288246
$data[11,16] = 5.667934548174888
$data[11,17] = 51.01141093357401
$data[11,18] = 0.2724179727278936
$data[11,19] = 0.2724179727278936

# Write the full A2:T13 block (header row 1 is unchanged) and let Excel resize the
# sheets used range / dimension automatically.
$ws.Range("A2:T13").Value = $data
